{"js": "// Update stack-trace line numbers / method references to reflect the\n// move from POI 3.17.0 to 4.0.1 (as captured in the regenerated JUnit\n// failure trace embedded in the document body).\n//\n// Each entry is an exact (old, new) text pair. Several edits are plain\n// single-line substitutions (changed line numbers); the last entry\n// replaces a large contiguous block of stack-trace lines (Maven\n// Surefire / Equinox launcher frames) with the new Eclipse JDT JUnit\n// runner frames.\nconst replacements = [\n  [\n    \"\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\",\n    \"\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)\"\n  ],\n  [\n    \"\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\",\n    \"\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)\"\n  ],\n  [\n    \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\",\n    \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)\"\n  ],\n  [\n    \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\",\n    \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)\"\n  ],\n  [\n    \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:333)\",\n    \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:345)\"\n  ],\n  [\n    \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:259)\",\n    \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:262)\"\n  ],\n  [\n    \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\",\n    \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)\"\n  ],\n  [\n    \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\",\n    \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)\"\n  ],\n  [\n    \"\\tat sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)\",\n    \"\\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\"\n  ],\n  [\n    \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\\n\" +\n    \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)\\n\" +\n    \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)\\n\" +\n    \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n    \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n    \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n    \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n    \"\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)\\n\" +\n    \"\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:156)\\n\" +\n    \"\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)\\n\" +\n    \"\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)\\n\" +\n    \"\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\" +\n    \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n    \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n    \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n    \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n    \"\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)\\n\" +\n    \"\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)\\n\" +\n    \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)\\n\" +\n    \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)\\n\" +\n    \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)\\n\" +\n    \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)\\n\" +\n    \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n    \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n    \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n    \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n    \"\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)\\n\" +\n    \"\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)\\n\" +\n    \"\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)\\n\" +\n    \"\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)\",\n\n    \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n    \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n    \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\" +\n    \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\" +\n    \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\" +\n    \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update stack-trace line numbers / method references to reflect the\n# move from POI 3.17.0 to 4.0.1 (as captured in the regenerated JUnit\n# failure trace embedded in the document body).\n\n$d = $word.ActiveDocument\n$tab = [char]9\n$nl = [char]10\n\nfunction Replace-Text($search, $replace) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $search\n    $find.Replacement.Text = $replace\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap (1=wdFindContinue), Format,\n    # ReplaceWith, Replace (2=wdReplaceAll)\n    $find.Execute($search, $false, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\n$s1 = $tab + \"at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\"\n$r1 = $tab + \"at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)\"\nReplace-Text $s1 $r1\n\n$s2 = $tab + \"at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\"\n$r2 = $tab + \"at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)\"\nReplace-Text $s2 $r2\n\n$s3 = $tab + \"at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\"\n$r3 = $tab + \"at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)\"\nReplace-Text $s3 $r3\n\n$s4 = $tab + \"at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\"\n$r4 = $tab + \"at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)\"\nReplace-Text $s4 $r4\n\n$s5 = $tab + \"at org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:333)\"\n$r5 = $tab + \"at org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:345)\"\nReplace-Text $s5 $r5\n\n$s6 = $tab + \"at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:259)\"\n$r6 = $tab + \"at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:262)\"\nReplace-Text $s6 $r6\n\n$s7 = $tab + \"at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\"\n$r7 = $tab + \"at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)\"\nReplace-Text $s7 $r7\n\n$s8 = $tab + \"at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\"\n$r8 = $tab + \"at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)\"\nReplace-Text $s8 $r8\n\n$s9 = $tab + \"at sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)\"\n$r9 = $tab + \"at sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\"\nReplace-Text $s9 $r9\n\n# Replace the large block of Maven Surefire / Equinox launcher frames\n# with the new Eclipse JDT JUnit runner frames.\n$removedBlock = (\n    $tab + \"at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\" + $nl +\n    $tab + \"at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)\" + $nl +\n    $tab + \"at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)\" + $nl +\n    $tab + \"at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\" + $nl +\n    $tab + \"at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\" + $nl +\n    $tab + \"at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\" + $nl +\n    $tab + \"at java.lang.reflect.Method.invoke(Method.java:498)\" + $nl +\n    $tab + \"at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)\" + $nl +\n    $tab + \"at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)\" + $nl +\n    $tab + \"at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)\" + $nl +\n    $tab + \"at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)\" + $nl +\n    $tab + \"at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\" + $nl +\n    $tab + \"at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\" + $nl +\n    $tab + \"at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\" + $nl +\n    $tab + \"at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\" + $nl +\n    $tab + \"at java.lang.reflect.Method.invoke(Method.java:498)\" + $nl +\n    $tab + \"at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)\" + $nl +\n    $tab + \"at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)\" + $nl +\n    $tab + \"at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)\" + $nl +\n    $tab + \"at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)\" + $nl +\n    $tab + \"at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)\" + $nl +\n    $tab + \"at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)\" + $nl +\n    $tab + \"at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\" + $nl +\n    $tab + \"at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\" + $nl +\n    $tab + \"at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\" + $nl +\n    $tab + \"at java.lang.reflect.Method.invoke(Method.java:498)\" + $nl +\n    $tab + \"at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)\" + $nl +\n    $tab + \"at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)\" + $nl +\n    $tab + \"at org.eclipse.equinox.launcher.Main.run(Main.java:1498)\" + $nl +\n    $tab + \"at org.eclipse.equinox.launcher.Main.main(Main.java:1471)\"\n)\n\n$addedBlock = (\n    $tab + \"at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\" + $nl +\n    $tab + \"at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\" + $nl +\n    $tab + \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\" + $nl +\n    $tab + \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\" + $nl +\n    $tab + \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\" + $nl +\n    $tab + \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\"\n)\n\nReplace-Text $removedBlock $addedBlock\n"}
